$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AddOpportunity")

# Duplicate formatting of row 2 into row 3 (same layout/style as existing data row)
$ws.Range("A2:AD2").Copy() | Out-Null
$ws.Range("A3:AD3").PasteSpecial(-4122) | Out-Null

# Now set the row 3 values (same as row 2 except Client, Subject, ClientOwnership,
# which describe the new test engagement: Pharmavite buying Food State via a Hedge Fund)
$ws.Range("A3").Value = "Pharmavite, LLC"
$ws.Range("B3").Value = "Food State, Inc."
$ws.Range("C3").Value = "Buyside"
$ws.Range("D3").Value = "BUS - Business Services"
$ws.Range("E3").Value = "Dealership & Rental Services"
$ws.Range("F3").Value = "No"
$ws.Range("G3").Value = "No"
$ws.Range("H3").Value = "Accountant"
$ws.Range("I3").Value = "No"
$ws.Range("J3").Value = "No"
$ws.Range("K3").Value = "AM"
$ws.Range("L3").Value = "HL Capital, Inc."
$ws.Range("M3").Value = "Do Not Disclose"
$ws.Range("N3").Value = "Mark Martin"
$ws.Range("O3").Value = "1000"
$ws.Range("P3").Value = "1000"
$ws.Range("Q3").Value = "1000"
$ws.Range("R3").Value = "Hedge Fund"
$ws.Range("S3").Value = "Family Office"
$ws.Range("T3").Value = "9999"
$ws.Range("U3").Value = "Test"
$ws.Range("V3").Value = "Chris Lord"
$ws.Range("W3").Value = "Yes, separate signed agreement"
$ws.Range("X3").Value = "Cleared"
$ws.Range("Y3").Value = "CF"
$ws.Range("Z3").Value = "Consulting"
$ws.Range("AA3").Value = "1000"
$ws.Range("AB3").Value = "1000"
$ws.Range("AC3").Value = "Mark Martin"
$ws.Range("AD3").Value = "Yes"

# New ClientOwnership value is wider than the column, so resize to fit its content
$ws.Columns.Item(18).AutoFit() | Out-Null

# Make AddOpportunity the active sheet/tab, with the new row's data selected
$ws.Activate()
$ws.Range("R9").Select() | Out-Null
